# Update "want to go" head-count column (F) for a handful of events.
# The workbook keeps the same event rows duplicated on the "全部类型"
# (All types) aggregate sheet, so every bump on 展览/演出 is mirrored there.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibition) sheet
$wsExhibit.Range("F19").Value = 870
$wsExhibit.Range("F22").Value = 323
$wsExhibit.Range("F28").Value = 2755
$wsExhibit.Range("F34").Value = 347
$wsExhibit.Range("F40").Value = 675

# 演出 (Show) sheet
$wsShow.Range("F2").Value = 33

# 全部类型 (All types) aggregate sheet - mirrors the rows above
$wsAll.Range("F6").Value  = 33
$wsAll.Range("F25").Value = 870
$wsAll.Range("F27").Value = 323
$wsAll.Range("F32").Value = 2755
$wsAll.Range("F38").Value = 347
$wsAll.Range("F43").Value = 675
